$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that looks numeric (e.g. "68.324.04", "1.00").
# Force them to stay text (matching the workbook author convention of inline/shared
# strings for every data cell) so Excel does not silently coerce them to numbers and
# drop significant trailing zeros / thousands-style separators.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.324.04"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.640.73"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.48"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.43"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.639.32"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.26"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.11"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000192"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.119.45"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.212.73"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.646.88"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.52"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  +2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.88"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.65"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.85"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.781.15"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "572.33"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.15"
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.43"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.05"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.38"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.373"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0334"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("E45").Value = "  +3.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.67"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.33"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.77"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.86"
